$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the position currently occupied by rows 253-254.
# This pushes the existing rows 253-308 down to 255-310 (matching the
# diff, which shows every existing row's data shifted down by two rows)
# and grows the sheet dimension from A1:T308 to A1:T310 automatically.
$ws.Rows("253:254").Insert()

# Populate the newly inserted row 253 with the weekly record added by
# this commit (Barraganete / Primera, week of 2022-01-07 -> serial 44543).
$ws.Cells.Item(253, 1).Value = 4
$ws.Cells.Item(253, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(253, 3).Value = "Los Lagos"
$ws.Cells.Item(253, 4).Value = 44543
$ws.Cells.Item(253, 5).Value = 10
$ws.Cells.Item(253, 6).Value = "Fruta"
$ws.Cells.Item(253, 7).Value = 100108
$ws.Cells.Item(253, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(253, 9).Value = 100108006
$ws.Cells.Item(253, 10).Value = "Plátano"
$ws.Cells.Item(253, 11).Value = "Barraganete"
$ws.Cells.Item(253, 12).Value = "Primera"
$ws.Cells.Item(253, 13).Value = 200
$ws.Cells.Item(253, 14).Value = 26000
$ws.Cells.Item(253, 15).Value = 27000
$ws.Cells.Item(253, 16).Value = 26500
$ws.Cells.Item(253, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(253, 18).Value = "Ecuador"
$ws.Cells.Item(253, 19).Value = 1325
$ws.Cells.Item(253, 20).Value = 20

# Populate the newly inserted row 254 with the second weekly record added
# by this commit (Sin especificar / Primera Pintón, same week).
$ws.Cells.Item(254, 1).Value = 4
$ws.Cells.Item(254, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(254, 3).Value = "Los Lagos"
$ws.Cells.Item(254, 4).Value = 44543
$ws.Cells.Item(254, 5).Value = 10
$ws.Cells.Item(254, 6).Value = "Fruta"
$ws.Cells.Item(254, 7).Value = 100108
$ws.Cells.Item(254, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(254, 9).Value = 100108006
$ws.Cells.Item(254, 10).Value = "Plátano"
$ws.Cells.Item(254, 11).Value = "Sin especificar"
$ws.Cells.Item(254, 12).Value = "Primera Pintón"
$ws.Cells.Item(254, 13).Value = 500
$ws.Cells.Item(254, 14).Value = 19000
$ws.Cells.Item(254, 15).Value = 20000
$ws.Cells.Item(254, 16).Value = 19500
$ws.Cells.Item(254, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(254, 18).Value = "Ecuador"
$ws.Cells.Item(254, 19).Value = 975
$ws.Cells.Item(254, 20).Value = 20
